# Cirrus-TA-Instant-form.docx : add the "800 GPU-h" alternative limit
# alongside the existing "Core-h" limit in four spots, per commit:
#   "access/instant.md,access/ta/Cirrus-TA-Instant-form.docx: added 800 GPUh limit"

$d = $word.ActiveDocument
$find = $d.Content.Find

# ---------------------------------------------------------------------
# 1) "You can apply for a maximum of 20000 Core-h."
#       -> "You can apply for a maximum of 20,000 Core-h or 800 GPU-h."
#    ("or" is italic)
# ---------------------------------------------------------------------
$found = $find.Execute("apply for a maximum of 20000 Core-h.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r = $find.Parent
    $matchStart = $r.Start
    $matchEnd = $r.End

    # "20" + "," + "000 ..."
    $prefixLen = "apply for a maximum of 20".Length
    $commaPos = $matchStart + $prefixLen
    $d.Range($commaPos, $commaPos).InsertBefore(",")
    $matchEnd = $matchEnd + 1

    # insert " or 800 GPU-h" right before the trailing period
    $periodPos = $matchEnd - 1
    $d.Range($periodPos, $periodPos).InsertBefore(" or 800 GPU-h")

    # italicise just the "or"
    $orStart = $periodPos + 1
    $orEnd = $orStart + 2
    $d.Range($orStart, $orEnd).Italic = 1
}

# ---------------------------------------------------------------------
# 2) "... receive up to 20000 Core-h to be used over a 6-month period ..."
#       -> "... receive up to 20,000 Core-h or 800 GPU-h to be used over a 6-month period ..."
# ---------------------------------------------------------------------
$found2 = $find.Execute("receive up to 20000 Core-h to be used", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2 = $find.Parent
    $matchStart2 = $r2.Start

    $prefixLen2 = "receive up to 20".Length
    $commaPos2 = $matchStart2 + $prefixLen2
    $d.Range($commaPos2, $commaPos2).InsertBefore(",")

    $beforeToLen2 = "receive up to 20,000 Core-h".Length
    $insertPos2 = $matchStart2 + $beforeToLen2
    $d.Range($insertPos2, $insertPos2).InsertBefore(" or 800 GPU-h")
}

# ---------------------------------------------------------------------
# 3) "[Enter total Core-h from resource calculation specified in Section 5 below]"
#       -> "[Enter total Core-h or GPU-h from resource calculation specified in Section 5 below]"
# ---------------------------------------------------------------------
$find.Execute("Core-h from resource calculation specified in Section 5", $false, $false, $false, $false, $false, $true, 1, $false, "Core-h or GPU-h from resource calculation specified in Section 5", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "... the total Core-h specified in 4.1 above) ..."
#       -> "... the total Core-h or GPU-h specified in 4.1 above) ..."
#    (this also absorbs/removes the stray "_GoBack" bookmark sitting in the run)
# ---------------------------------------------------------------------
$find.Execute("Core-h specified in 4.1 above", $false, $false, $false, $false, $false, $true, 1, $false, "Core-h or GPU-h specified in 4.1 above", 2) | Out-Null

# ---------------------------------------------------------------------
# Cosmetic cleanups present in the target revision (no visible text change,
# just run / spell-check-marker tidy-up caused by Word re-saving the file).
# Re-applying the unchanged text via Find/Replace reflows the runs the same
# way a normal Word edit pass would.
# ---------------------------------------------------------------------
$find.Execute("By default, projects are granted 500 GiB of disk space. If you require more than this you should specify this here and justify t", $false, $false, $false, $false, $false, $true, 1, $false, "By default, projects are granted 500 GiB of disk space. If you require more than this you should specify this here and justify t", 2) | Out-Null

$find.Execute("[Specify disk space requirements if larger than 500 GiB]", $false, $false, $false, $false, $false, $true, 1, $false, "[Specify disk space requirements if larger than 500 GiB]", 2) | Out-Null

$find.Execute(" if more than 500 GiB", $false, $false, $false, $false, $false, $true, 1, $false, " if more than 500 GiB", 2) | Out-Null

$find.Execute("[Enter estimated total size in MiB/GiB/TiB]", $false, $false, $false, $false, $false, $true, 1, $false, "[Enter estimated total size in MiB/GiB/TiB]", 2) | Out-Null
$find.Execute("[Enter estimated total size in MiB/GiB/TiB]", $false, $false, $false, $false, $false, $true, 1, $false, "[Enter estimated total size in MiB/GiB/TiB]", 2) | Out-Null

Write-Output "done"
